$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet, matching the
# position of "2017-05-02" in the target workbook (it must come last).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2017-05-02"

# Header row
$ws.Range("A1").Value = "CreatedBy"
$ws.Range("B1").Value = "DataCompleted"
$ws.Range("C1").Value = "DataIncomplete"
$ws.Range("D1").Value = "Total"

$data = @(
    ,@('Dr Faeiz', 1, 0, 1)
    ,@('delina', 2, 0, 2)
    ,@('Aziani', 3, 0, 3)
    ,@('Jennifer Kaur', 6, 0, 6)
    ,@('Faeiz', 11, 0, 11)
    ,@('Ling Kuok Wei', 15, 0, 15)
    ,@('Dr Masliyana', 21, 0, 21)
    ,@('Hui Che', 21, 0, 21)
    ,@('Dr Faeiz Syezri Adzmin bin Jaaffar', 21, 1, 22)
    ,@('Afiq Firdaus', 23, 1, 24)
    ,@('Fadzli', 25, 0, 25)
    ,@('Ang SH', 26, 0, 26)
    ,@('Suhayl', 28, 0, 28)
    ,@('Aimi Nadiah Jamel', 27, 2, 29)
    ,@('Izzat', 29, 0, 29)
    ,@('Hui Yi', 37, 0, 37)
    ,@('Adlan', 38, 1, 39)
    ,@('Annas', 39, 0, 39)
    ,@('Eliza', 40, 0, 40)
    ,@('Munirah', 40, 0, 40)
    ,@('Noor hidayah', 40, 0, 40)
    ,@('Siti Aminah', 38, 2, 40)
    ,@('Izzati', 43, 0, 43)
    ,@('Musfirah', 43, 0, 43)
    ,@('Noor Amalina', 44, 0, 44)
    ,@('Natrah', 49, 0, 49)
    ,@('Nursyuhaida', 49, 0, 49)
    ,@('Hadi', 50, 0, 50)
    ,@('Aisyah', 50, 1, 51)
    ,@('Aizat', 50, 1, 51)
    ,@('Delina', 51, 0, 51)
    ,@('Danial', 52, 0, 52)
    ,@('Philip', 54, 0, 54)
    ,@('Dr Richard', 55, 0, 55)
    ,@('Helmi', 59, 0, 59)
    ,@('Yhyviyaa', 60, 0, 60)
    ,@('Amira', 62, 0, 62)
    ,@('Syahirah', 70, 0, 70)
    ,@('Hooi Fan', 72, 0, 72)
    ,@('Mardhiah', 74, 0, 74)
    ,@('Michelle', 75, 0, 75)
    ,@('Nabilah Iffah', 75, 0, 75)
    ,@('Yvonne', 74, 1, 75)
    ,@('Thivashini', 77, 0, 77)
    ,@('Aminiril anisah', 79, 0, 79)
    ,@('Geetha Krishnan', 84, 0, 84)
    ,@('Tan khai shin', 84, 0, 84)
    ,@('Aishah', 86, 0, 86)
    ,@('Pui Yee', 88, 0, 88)
    ,@('Jocelyn', 94, 0, 94)
    ,@('Sree Durga', 94, 0, 94)
    ,@('Hoong Ping', 113, 0, 113)
    ,@('Nurjannah', 111, 2, 113)
    ,@('Jia yi', 122, 0, 122)
    ,@('Yi Shin', 145, 0, 145)
    ,@('Jacelyn', 230, 0, 230)
    ,@('Total', 3149, 12, 3161)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}
